$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.932899999999994
$ws.Range("B7").Value = 4.7822
$ws.Range("A8").Value = -22.28000000000002
$ws.Range("A10").Value = -21.8266
$ws.Range("E10").Value = 15.76279999999999
$ws.Range("A12").Value = -21.55110000000001
$ws.Range("E12").Value = 17.70100000000001
$ws.Range("E13").Value = 16.60070000000001
$ws.Range("E14").Value = 16.71860000000001
$ws.Range("B15").Value = 4.430399999999997
$ws.Range("A18").Value = -22.5187
$ws.Range("B18").Value = 4.479099999999996
$ws.Range("D18").Value = -8.178599999999994
$ws.Range("D19").Value = -9.084499999999986
$ws.Range("B20").Value = 9.795599999999986
$ws.Range("D27").Value = -8.862199999999996
$ws.Range("B29").Value = 5.261900000000002
$ws.Range("E29").Value = 17.19
$ws.Range("B30").Value = 4.581300000000001
$ws.Range("B31").Value = 5.875800000000003
$ws.Range("D31").Value = -8.270599999999995
$ws.Range("E32").Value = 16.0622
$ws.Range("E35").Value = 16.5304
$ws.Range("A37").Value = -20.26079999999999
$ws.Range("D38").Value = -8.565700000000007
$ws.Range("B40").Value = 8.826599999999996
$ws.Range("D42").Value = -8.853799999999996
$ws.Range("E43").Value = 17.4201
$ws.Range("D44").Value = -7.8283
$ws.Range("D47").Value = -7.6607
$ws.Range("E48").Value = 17.41540000000002
$ws.Range("E49").Value = 15.7463
$ws.Range("B50").Value = 4.723200000000002
$ws.Range("E50").Value = 16.462
$ws.Range("A55").Value = -21.7654
$ws.Range("E56").Value = 16.0635
$ws.Range("D58").Value = -8.206
$ws.Range("D65").Value = -7.712099999999998
$ws.Range("A68").Value = -21.51700000000001
$ws.Range("B68").Value = 4.7992
$ws.Range("E69").Value = 17.39040000000002
$ws.Range("D73").Value = -7.825399999999997
$ws.Range("B76").Value = 6.6085
$ws.Range("A77").Value = -20.88079999999999
$ws.Range("A78").Value = -20.16289999999998
$ws.Range("A81").Value = -22.02130000000001
$ws.Range("E81").Value = 16.57340000000001
$ws.Range("A82").Value = -21.7292
$ws.Range("B87").Value = 4.701699999999994
$ws.Range("B88").Value = 4.777899999999997
$ws.Range("D90").Value = -8.094000000000003
$ws.Range("E92").Value = 18.57470000000002
$ws.Range("D94").Value = -6.824199999999999
$ws.Range("D95").Value = -7.942899999999998
$ws.Range("B96").Value = 4.820700000000007
$ws.Range("B98").Value = 6.047899999999998
$ws.Range("B101").Value = 9.023299999999999
$ws.Range("D101").Value = -7.541399999999999
$ws.Range("B102").Value = 8.7745
